$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 6 data
$ws.Range("A6").Value = 6
$ws.Range("B6").Value = "Thông báo qua xưởng làm"
$ws.Range("C6").Value = '<p style="text-align: center;"><span style="color: rgb(255, 0, 0);">Yêu cầu 100% qua xưởng</span></p>'
$ws.Range("E6").Value = "Ban Điều hành"
$ws.Range("F6").Value = "14/09/2022 04:21"
$ws.Range("G6").Value = "https://drive.google.com/drive/u/0/my-drive"
